# Applies the cryptocurrency price/volume update described in the commit.
# Values are written as literal text (matching the original inline-string cells),
# including numeric-looking Price figures, so formatting such as trailing zeros
# (e.g. "116.40") and scientific-looking decimals (e.g. "0.0000110") is preserved
# exactly rather than being auto-converted to numbers by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.760.73'
$ws.Range("E2").Value = '  -0.34%  '
# Row 3
$ws.Range("D3").Value = '2.288.59'
$ws.Range("E3").Value = '  -0.83%  '
# Row 4
$ws.Range("E4").Value = '  +0.13%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '116.40'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +13.64%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '269.34'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.70%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.627'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.06%  '
# Row 8
$ws.Range("E8").Value = '  +0.09%  '
# Row 9
$ws.Range("E9").Value = '  +1.75%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '49.09'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.96%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0945'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.95%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.98'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +13.17%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.107'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.18%  '
# Row 14
$ws.Range("E14").Value = '  -0.72%  '
# Row 15
$ws.Range("D15").Value = '2.631.73'
$ws.Range("E15").Value = '  -0.64%  '
# Row 16
$ws.Range("E16").Value = '  +1.38%  '
# Row 17
$ws.Range("D17").Value = '2.282.79'
$ws.Range("E17").Value = '  -0.45%  '
# Row 18
$ws.Range("D18").Value = '43.604.38'
$ws.Range("E18").Value = '  -0.47%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000110'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.10%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.98'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +11.54%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.54'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.18%  '
# Row 22
$ws.Range("E22").Value = '  -1.56%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.14'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +10.72%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '233.53'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.06%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.96'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.17%  '
# Row 26
$ws.Range("E26").Value = '  -0.06%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.69'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.97%  '
# Row 28
$ws.Range("E28").Value = '  +3.79%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '42.03'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +9.59%  '
# Row 30
$ws.Range("E30").Value = '  -2.16%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.24'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.06%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '173.69'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.11%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0940'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.97%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '21.60'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.31%  '
# Row 35
$ws.Range("E35").Value = '  +4.72%  '
# Row 36
$ws.Range("E36").Value = '  +0.18%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.75'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.58%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0359'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.81%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.108'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.00%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.86'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.91%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.57'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +19.08%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.86'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +15.55%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.44'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.62%  '
# Row 44
$ws.Range("E44").Value = '  +2.29%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.38'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +21.59%  '
# Row 46
$ws.Range("E46").Value = '  +0.06%  '
# Row 47
$ws.Range("E47").Value = '  +0.37%  '
# Row 48
$ws.Range("E48").Value = '  -1.18%  '
# Row 49
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '102.81'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.04%  '
# Row 50
$ws.Range("B50").Value = 'TrustWalletToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.26'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.66%  '
# Row 51
$ws.Range("E51").Value = '  -1.87%  '
